$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a weekly price log for "Coliflor" at the La Palmera de La
# Serena terminal. Each reporting date occupies two consecutive rows (one
# "Primera" quality row followed by one "Segunda" quality row). This commit
# adds one more week of data (two new rows) right above the existing block
# that starts at row 346, and every following row shifts down by two.

# 1) Insert two blank rows at row 346 - this pushes the old rows 346.. down
#    to 348.. and grows the sheet from 462 to 464 used rows.
$ws.Rows.Item(346).Resize(2).Insert()

# 2) The two new blank rows need the same formatting / constant columns
#    (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
#    Calidad, Unidad de comercializacion, Origen, Kg o Unidades,
#    Clasificacion) as the rows directly below them (which now hold what used
#    to be row 346 and row 347, i.e. the same "Primera"/"Segunda" pair
#    pattern). Copy those rows down into the new blank ones first, then we
#    will overwrite just the cells that actually carry new data.
$ws.Rows.Item(348).Copy()
$ws.Rows.Item(346).PasteSpecial()
$ws.Rows.Item(349).Copy()
$ws.Rows.Item(347).PasteSpecial()

# 3) Now set the new week's actual reported figures for the two new rows.
# Row 346 - "Primera" quality
$ws.Range("D346").Value = 44524
$ws.Range("J346").Value = 3340
$ws.Range("K346").Value = 600
$ws.Range("L346").Value = 700
$ws.Range("M346").Value = 650
$ws.Range("P346").Value = 650

# Row 347 - "Segunda" quality
$ws.Range("D347").Value = 44524
$ws.Range("J347").Value = 1600
$ws.Range("K347").Value = 500
$ws.Range("L347").Value = 550
$ws.Range("M347").Value = 525
$ws.Range("P347").Value = 525

$excel.CutCopyMode = 0
